$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7
$ws.Range("AH2").Value = 9
$ws.Range("AI2").Value = 19
$ws.Range("AL2").Value = 8
$ws.Range("AM2").Value = 8
$ws.Range("AN2").Value = 9
$ws.Range("G5").Value = 2.05
$ws.Range("I5").Value = 4.1
$ws.Range("T5").Value = 1.16
$ws.Range("AR5").Value = 2.25
$ws.Range("AS5").Value = 1.65
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("S9").Value = 3
$ws.Range("T9").Value = 1.38
$ws.Range("AR9").Value = 1.57
$ws.Range("AS9").Value = 2.36
$ws.Range("G10").Value = 2.7
$ws.Range("I10").Value = 2.8
$ws.Range("J10").Value = 3.5
$ws.Range("K10").Value = 1.91
$ws.Range("L10").Value = 3.6
$ws.Range("O10").Value = 1.5
$ws.Range("P10").Value = 2.5
$ws.Range("AA10").Value = 7
$ws.Range("AB10").Value = 12
$ws.Range("AC10").Value = 11
$ws.Range("AD10").Value = 29
$ws.Range("AE10").Value = 26
$ws.Range("AL10").Value = 7
$ws.Range("AM10").Value = 12
$ws.Range("AN10").Value = 11
$ws.Range("AO10").Value = 29
$ws.Range("AP10").Value = 26
$ws.Range("AR10").Value = 1.9
$ws.Range("AS10").Value = 1.95
$ws.Range("G11").Value = 2.87
$ws.Range("H11").Value = 2.67
$ws.Range("I11").Value = 2.77
$ws.Range("J11").Value = 3.65
$ws.Range("K11").Value = 1.82
$ws.Range("L11").Value = 3.5
$ws.Range("M11").Value = 1.15
$ws.Range("N11").Value = 4.75
$ws.Range("O11").Value = 1.62
$ws.Range("P11").Value = 2.15
$ws.Range("Q11").Value = 2.82
$ws.Range("R11").Value = 1.38
$ws.Range("U11").Value = 5.2
$ws.Range("V11").Value = 1.13
$ws.Range("W11").Value = 1.62
$ws.Range("X11").Value = 2.15
$ws.Range("Y11").Value = 2.2
$ws.Range("Z11").Value = 1.6
$ws.Range("AA11").Value = 6.1
$ws.Range("AB11").Value = 12.5
$ws.Range("AC11").Value = 11.5
$ws.Range("AD11").Value = 37
$ws.Range("AE11").Value = 35
$ws.Range("AF11").Value = 60
$ws.Range("AG11").Value = 4.75
$ws.Range("AH11").Value = 5.5
$ws.Range("AI11").Value = 19
$ws.Range("AJ11").Value = 150
$ws.Range("AK11").Value = 101
$ws.Range("AL11").Value = 6.1
$ws.Range("AM11").Value = 12.5
$ws.Range("AN11").Value = 11
$ws.Range("AO11").Value = 35
$ws.Range("AP11").Value = 32
$ws.Range("AQ11").Value = 55
$ws.Range("G12").Value = 3.35
$ws.Range("H12").Value = 2.77
$ws.Range("I12").Value = 2.4
$ws.Range("J12").Value = 3.9
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 5.5
$ws.Range("O12").Value = 1.47
$ws.Range("P12").Value = 2.52
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 1.53
$ws.Range("U12").Value = 4.1
$ws.Range("V12").Value = 1.2
$ws.Range("Y12").Value = 1.91
$ws.Range("Z12").Value = 1.8
$ws.Range("AA12").Value = 7.9
$ws.Range("AB12").Value = 16.5
$ws.Range("AC12").Value = 11.5
$ws.Range("AE12").Value = 35
$ws.Range("AF12").Value = 45
$ws.Range("AG12").Value = 5.5
$ws.Range("AH12").Value = 5.4
$ws.Range("AI12").Value = 15
$ws.Range("AK12").Value = 800
$ws.Range("AM12").Value = 11
$ws.Range("AO12").Value = 26
$ws.Range("G13").Value = 2.07
$ws.Range("H13").Value = 2.82
$ws.Range("I13").Value = 4.1
$ws.Range("J13").Value = 2.82
$ws.Range("K13").Value = 1.85
$ws.Range("L13").Value = 4.65
$ws.Range("M13").Value = 1.14
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 1.57
$ws.Range("P13").Value = 2.25
$ws.Range("Q13").Value = 2.67
$ws.Range("R13").Value = 1.42
$ws.Range("U13").Value = 4.9
$ws.Range("V13").Value = 1.14
$ws.Range("W13").Value = 1.6
$ws.Range("X13").Value = 2.2
$ws.Range("Y13").Value = 2.22
$ws.Range("Z13").Value = 1.6
$ws.Range("AA13").Value = 5
$ws.Range("AB13").Value = 8.25
$ws.Range("AC13").Value = 9.5
$ws.Range("AD13").Value = 19.5
$ws.Range("AE13").Value = 23
$ws.Range("AF13").Value = 50
$ws.Range("AG13").Value = 5
$ws.Range("AH13").Value = 5.7
$ws.Range("AI13").Value = 19.5
$ws.Range("AJ13").Value = 150
$ws.Range("AK13").Value = 101
$ws.Range("AL13").Value = 8.25
$ws.Range("AM13").Value = 21
$ws.Range("AN13").Value = 14.5
$ws.Range("AO13").Value = 70
$ws.Range("AP13").Value = 50
$ws.Range("AQ13").Value = 65
$ws.Range("G15").Value = 1.7
$ws.Range("H15").Value = 3.6
$ws.Range("O16").Value = 1.4
$ws.Range("P16").Value = 3
$ws.Range("Q16").Value = 2.2
$ws.Range("R16").Value = 1.67
$ws.Range("G17").Value = 3.5
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 1.08
$ws.Range("N17").Value = 7.5
$ws.Range("O17").Value = 1.44
$ws.Range("P17").Value = 2.75
$ws.Range("Q17").Value = 2.38
$ws.Range("R17").Value = 1.57
$ws.Range("U17").Value = 4.5
$ws.Range("V17").Value = 1.2
$ws.Range("W17").Value = 1.53
$ws.Range("X17").Value = 2.38
$ws.Range("Y17").Value = 2
$ws.Range("Z17").Value = 1.73
$ws.Range("AG17").Value = 7.5
$ws.Range("AI17").Value = 17
$ws.Range("AL17").Value = 6.5
$ws.Range("AN17").Value = 9.5
$ws.Range("AQ17").Value = 34
$ws.Range("AR17").Value = 1.8
$ws.Range("AS17").Value = 2.05
$ws.Range("G19").Value = 2.5
$ws.Range("I19").Value = 2.75
$ws.Range("L19").Value = 3.5
$ws.Range("N19").Value = 9.5
$ws.Range("Q19").Value = 2.1
$ws.Range("R19").Value = 1.73
$ws.Range("AD19").Value = 23
$ws.Range("AO19").Value = 29
$ws.Range("AP19").Value = 23
$ws.Range("M23").Value = 1.08
$ws.Range("N23").Value = 8
$ws.Range("Q23").Value = 2.4
$ws.Range("R23").Value = 1.5
$ws.Range("AR23").Value = 1.83
$ws.Range("AS23").Value = 1.98
$ws.Range("G25").Value = 2.63
$ws.Range("I25").Value = 2.7
$ws.Range("J25").Value = 3.25
$ws.Range("K25").Value = 2.1
$ws.Range("L25").Value = 3.4
$ws.Range("M25").Value = 1.06
$ws.Range("N25").Value = 10
$ws.Range("O25").Value = 1.3
$ws.Range("P25").Value = 3.4
$ws.Range("Q25").Value = 2.03
$ws.Range("R25").Value = 1.78
$ws.Range("U25").Value = 3.5
$ws.Range("V25").Value = 1.29
$ws.Range("Y25").Value = 1.73
$ws.Range("Z25").Value = 2
$ws.Range("AC25").Value = 10
$ws.Range("AE25").Value = 21
$ws.Range("AF25").Value = 29
$ws.Range("AG25").Value = 9
$ws.Range("AI25").Value = 13
$ws.Range("AJ25").Value = 41
$ws.Range("AK25").Value = 201
$ws.Range("AL25").Value = 9
$ws.Range("AN25").Value = 11
$ws.Range("AQ25").Value = 29
